# [Anmol Singh] Interface implemented successfully
# Refresh the stored DOB timestamp (same calendar date, new time-of-day
# fraction) and rotate the stored password hash for every customer /
# driver / owner record.

$wb = $excel.ActiveWorkbook

$newDob = 28430.95315505787

# --- Customer sheet -------------------------------------------------
$ws = $wb.Worksheets.Item("Customer")

$ws.Range("E2").Value2 = $newDob
$ws.Range("E3").Value2 = $newDob
$ws.Range("E4").Value2 = $newDob
$ws.Range("E5").Value2 = $newDob
$ws.Range("E6").Value2 = $newDob
$ws.Range("E7").Value2 = $newDob
$ws.Range("E8").Value2 = $newDob

$ws.Range("I2").Value2 = "l+BDBpDhP3PWK4Fdk1V0gkk7jaTT7MYGlNeda1DFmjs="
$ws.Range("I3").Value2 = "JHR8RykEfHhWuFvr2qFg7UOH0poHNV6epBwsNfMvbos="
$ws.Range("I4").Value2 = "Nc2LEU0rxCIZcxOFSLpqvpvrbnIkyc00SP+VtvQBB3k="
$ws.Range("I5").Value2 = "sl3p02JrjNPZgMMxG3UeO0JJN2tkTf7kpiVVewyoxbU="
$ws.Range("I6").Value2 = "B115B8/yKg64ryK4ohFmNI6rR3hGcmHJCR6jWRG2b5k="
$ws.Range("I7").Value2 = "2ExNQn0Sq8qtQa7hbNPPx9yqg8t56+7xl0FzmHdf8PE="
$ws.Range("I8").Value2 = "V8bBwkyyhFfzJr4tq55HNNxe7usmX7P1VUF7OkRe+TM="

# --- Driver sheet -----------------------------------------------------
$ws = $wb.Worksheets.Item("Driver")

$ws.Range("E2").Value2 = $newDob
$ws.Range("E3").Value2 = $newDob
$ws.Range("E4").Value2 = $newDob
$ws.Range("E5").Value2 = $newDob
$ws.Range("E6").Value2 = $newDob

$ws.Range("J2").Value2 = "xvPROODH23AWNs6tXxgQ98YaMoqyYdjCFo5+YUSrpfs="
$ws.Range("J3").Value2 = "xPCZJISMu3W5wCYgS4zrOB1nY03UzSGYcJcyNhg9o2w="
$ws.Range("J4").Value2 = "erlEc/fLo3OLNwoVdIfqtocMPMSjU37j5M5Kvpg8h6g="
$ws.Range("J5").Value2 = "TO2sBnBGPCnXowSmxLH/nrmKaIMkP5Wkd+8JRP1v2+s="
$ws.Range("J6").Value2 = "caumWCZ3K+hwiBDrN7Q0P+9yL/qKgfNRqFO8fenzv2o="

# --- Owner sheet ------------------------------------------------------
$ws = $wb.Worksheets.Item("Owner")

$ws.Range("E2").Value2 = $newDob
$ws.Range("E3").Value2 = $newDob
$ws.Range("E4").Value2 = $newDob

$ws.Range("J2").Value2 = "F9iYinmeDhDQSVn6dxrCoVUscg122nS/fDmYP0DW25I="
$ws.Range("J3").Value2 = "RcQlz8KhXIENPczp0emYjGL2l3o5h0YD48K9WNZuEWI="
$ws.Range("J4").Value2 = "AisijCSPzx3NL3S7Cmw4tcKGtoAIKIBJ+3hpbL7GAsU="
